# Apply updated cryptocurrency price/volume data to the "cryptos" sheet.
# Values are written as literal text (matching the workbook's existing
# inline-string cell convention) rather than being auto-coerced into
# numbers, so formatting such as "57.366.82" or "  -4.29%  " is preserved
# exactly as scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell that is already known to carry the sheet's default (unstyled)
# text format - used to restore styling after the NumberFormat="@" trick
# below, so we don't leave stray "Text" number formats behind.
$defaultStyle = $ws.Range("B2").Style

$updates = @(
    @{ Ref = "D2"; Value = "57.366.82" },
    @{ Ref = "E2"; Value = "  -4.29%  " },
    @{ Ref = "D3"; Value = "2.918.36" },
    @{ Ref = "E3"; Value = "  -2.39%  " },
    @{ Ref = "E4"; Value = "  +0.07%  " },
    @{ Ref = "D5"; Value = "546.75" },
    @{ Ref = "E5"; Value = "  -4.18%  " },
    @{ Ref = "D6"; Value = "128.89" },
    @{ Ref = "E6"; Value = "  +2.98%  " },
    @{ Ref = "E7"; Value = "  +0.11%  " },
    @{ Ref = "D8"; Value = "0.509" },
    @{ Ref = "E8"; Value = "  +1.56%  " },
    @{ Ref = "D9"; Value = "2.911.96" },
    @{ Ref = "E9"; Value = "  -2.53%  " },
    @{ Ref = "E10"; Value = "  -3.58%  " },
    @{ Ref = "E11"; Value = "  -6.32%  " },
    @{ Ref = "D12"; Value = "0.445" },
    @{ Ref = "E12"; Value = "  +1.27%  " },
    @{ Ref = "E13"; Value = "  -0.21%  " },
    @{ Ref = "D14"; Value = "32.69" },
    @{ Ref = "E14"; Value = "  +0.68%  " },
    @{ Ref = "E15"; Value = "  -0.01%  " },
    @{ Ref = "D16"; Value = "3.402.75" },
    @{ Ref = "E16"; Value = "  -2.14%  " },
    @{ Ref = "E17"; Value = "  +5.48%  " },
    @{ Ref = "D18"; Value = "2.910.51" },
    @{ Ref = "E18"; Value = "  -2.32%  " },
    @{ Ref = "D19"; Value = "57.414.02" },
    @{ Ref = "E19"; Value = "  -4.24%  " },
    @{ Ref = "D20"; Value = "415.71" },
    @{ Ref = "E20"; Value = "  -2.18%  " },
    @{ Ref = "E21"; Value = "  +0.04%  " },
    @{ Ref = "D22"; Value = "0.680" },
    @{ Ref = "E22"; Value = "  +2.07%  " },
    @{ Ref = "E23"; Value = "  -1.37%  " },
    @{ Ref = "D24"; Value = "12.96" },
    @{ Ref = "E24"; Value = "  +0.08%  " },
    @{ Ref = "D25"; Value = "79.57" },
    @{ Ref = "E25"; Value = "  +0.56%  " },
    @{ Ref = "D26"; Value = "1.00" },
    @{ Ref = "E26"; Value = "  -0.08%  " },
    @{ Ref = "D27"; Value = "0.999" },
    @{ Ref = "E27"; Value = "  +0.03%  " },
    @{ Ref = "D28"; Value = "2.47" },
    @{ Ref = "E28"; Value = "  -2.37%  " },
    @{ Ref = "E29"; Value = "  +2.00%  " },
    @{ Ref = "E30"; Value = "  +1.22%  " },
    @{ Ref = "D31"; Value = "25.13" },
    @{ Ref = "E31"; Value = "  +0.06%  " },
    @{ Ref = "D32"; Value = "5.92" },
    @{ Ref = "E32"; Value = "  -4.33%  " },
    @{ Ref = "D33"; Value = "0.0964" },
    @{ Ref = "E33"; Value = "  +2.80%  " },
    @{ Ref = "D34"; Value = "5.62" },
    @{ Ref = "E34"; Value = "  +0.22%  " },
    @{ Ref = "E35"; Value = "  +0.27%  " },
    @{ Ref = "D36"; Value = "2.05" },
    @{ Ref = "E36"; Value = "  -0.23%  " },
    @{ Ref = "D37"; Value = "47.74" },
    @{ Ref = "E37"; Value = "  -4.48%  " },
    @{ Ref = "E38"; Value = "  +4.06%  " },
    @{ Ref = "E39"; Value = "  +2.08%  " },
    @{ Ref = "D40"; Value = "2.53" },
    @{ Ref = "E40"; Value = "  +2.99%  " },
    @{ Ref = "B41"; Value = "Bittensor" },
    @{ Ref = "C41"; Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao" },
    @{ Ref = "D41"; Value = "373.77" },
    @{ Ref = "E41"; Value = "  -1.26%  " },
    @{ Ref = "B42"; Value = "Kaspa" },
    @{ Ref = "C42"; Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas" },
    @{ Ref = "D42"; Value = "0.106" },
    @{ Ref = "E42"; Value = "  -1.70%  " },
    @{ Ref = "B43"; Value = "VeChain" },
    @{ Ref = "C43"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" },
    @{ Ref = "D43"; Value = "0.0343" },
    @{ Ref = "E43"; Value = "  -3.10%  " },
    @{ Ref = "D44"; Value = "2.655.33" },
    @{ Ref = "E44"; Value = "  -0.37%  " },
    @{ Ref = "D46"; Value = "122.06" },
    @{ Ref = "E46"; Value = "  +1.63%  " },
    @{ Ref = "D47"; Value = "0.236" },
    @{ Ref = "E47"; Value = "  +0.85%  " },
    @{ Ref = "D48"; Value = "0.108" },
    @{ Ref = "E48"; Value = "  +1.41%  " },
    @{ Ref = "D49"; Value = "1.95" },
    @{ Ref = "E49"; Value = "  -2.19%  " },
    @{ Ref = "D50"; Value = "23.04" },
    @{ Ref = "E50"; Value = "  -2.50%  " },
    @{ Ref = "E51"; Value = "  -0.50%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    # Force text storage so numeric-looking strings (prices, percentages)
    # aren't reinterpreted/rounded as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = $defaultStyle
}
